$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column keeps text formatting so values like "1.00" or
# "0.140" are not silently coerced to numbers (dropping formatting).
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range('D2').Value = '67.232.37'
$ws.Range('E2').Value = '  -3.04%  '
$ws.Range('D3').Value = '3.538.62'
$ws.Range('E3').Value = '  -3.63%  '
$ws.Range('E4').Value = '  +0.02%  '
$ws.Range('D5').Value = '609.96'
$ws.Range('E5').Value = '  -5.58%  '
$ws.Range('D6').Value = '153.77'
$ws.Range('E6').Value = '  -2.72%  '
$ws.Range('D7').Value = '3.536.55'
$ws.Range('E7').Value = '  -3.62%  '
$ws.Range('E8').Value = '  +0.03%  '
$ws.Range('E9').Value = '  -2.70%  '
$ws.Range('D10').Value = '0.140'
$ws.Range('E10').Value = '  -2.70%  '
$ws.Range('D11').Value = '6.84'
$ws.Range('E11').Value = '  -3.49%  '
$ws.Range('E12').Value = '  -3.29%  '
$ws.Range('D13').Value = '0.0000222'
$ws.Range('E13').Value = '  -3.36%  '
$ws.Range('D14').Value = '4.139.72'
$ws.Range('E14').Value = '  -3.52%  '
$ws.Range('D15').Value = '31.88'
$ws.Range('E15').Value = '  -1.66%  '
$ws.Range('D16').Value = '3.535.22'
$ws.Range('E16').Value = '  -3.37%  '
$ws.Range('D17').Value = '67.174.41'
$ws.Range('E17').Value = '  -3.11%  '
$ws.Range('E18').Value = '  +1.14%  '
$ws.Range('D19').Value = '6.33'
$ws.Range('E19').Value = '  -1.71%  '
$ws.Range('D20').Value = '15.41'
$ws.Range('E20').Value = '  -2.49%  '
$ws.Range('D21').Value = '446.00'
$ws.Range('E21').Value = '  -4.15%  '
$ws.Range('D22').Value = '9.27'
$ws.Range('E22').Value = '  -6.92%  '
$ws.Range('D23').Value = '0.631'
$ws.Range('E23').Value = '  -1.87%  '
$ws.Range('D24').Value = '77.99'
$ws.Range('E24').Value = '  -1.62%  '
$ws.Range('D25').Value = '3.679.17'
$ws.Range('E25').Value = '  -3.64%  '
$ws.Range('E26').Value = '  +0.06%  '
$ws.Range('D27').Value = '0.0000122'
$ws.Range('E27').Value = '  -1.33%  '
$ws.Range('D28').Value = '10.24'
$ws.Range('E28').Value = '  -4.24%  '
$ws.Range('D29').Value = '8.27'
$ws.Range('E29').Value = '  -7.43%  '
$ws.Range('E30').Value = '  -3.27%  '
$ws.Range('E31').Value = '  -0.39%  '
$ws.Range('D32').Value = '1.00'
$ws.Range('E32').Value = '  -0.03%  '
$ws.Range('D33').Value = '25.77'
$ws.Range('E33').Value = '  -4.03%  '
$ws.Range('D34').Value = '0.159'
$ws.Range('E34').Value = '  -1.52%  '
$ws.Range('B35').Value = 'NEARProtocol'
$ws.Range('C35').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D35').Value = '6.17'
$ws.Range('E35').Value = '  -3.18%  '
$ws.Range('B36').Value = 'ImmutableX'
$ws.Range('C36').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D36').Value = '1.87'
$ws.Range('E36').Value = '  -5.60%  '
$ws.Range('D37').Value = '3.532.47'
$ws.Range('E37').Value = '  -3.61%  '
$ws.Range('D38').Value = '8.03'
$ws.Range('E38').Value = '  -3.90%  '
$ws.Range('E39').Value = '  +0.04%  '
$ws.Range('D40').Value = '1.00'
$ws.Range('E40').Value = '  -0.01%  '
$ws.Range('D41').Value = '175.42'
$ws.Range('E41').Value = '  -1.82%  '
$ws.Range('E42').Value = '  -2.48%  '
$ws.Range('D43').Value = '5.57'
$ws.Range('E43').Value = '  -4.66%  '
$ws.Range('D44').Value = '0.0865'
$ws.Range('E44').Value = '  -2.83%  '
$ws.Range('E45').Value = '  -3.37%  '
$ws.Range('D46').Value = '45.69'
$ws.Range('E46').Value = '  -2.40%  '
$ws.Range('D47').Value = '27.58'
$ws.Range('E47').Value = '  -1.57%  '
$ws.Range('D48').Value = '2.62'
$ws.Range('E48').Value = '  -2.45%  '
$ws.Range('D49').Value = '1.22'
$ws.Range('E49').Value = '  -1.16%  '
$ws.Range('B50').Value = 'SuiNetwork'
$ws.Range('C50').Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range('D50').Value = '1.03'
$ws.Range('E50').Value = '  -2.35%  '
$ws.Range('B51').Value = 'Cosmos'
$ws.Range('C51').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D51').Value = '7.58'
$ws.Range('E51').Value = '  -2.45%  '
